# Insert a new data row right after the existing "Vega Modelo de Temuco - Albahaca"
# row for 2021-08-12 (row 64), shifting every subsequent row down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(64).Insert()

$ws.Cells.Item(64, 1).Value = 10
$ws.Cells.Item(64, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(64, 3).Value = "La Araucanía"
$ws.Cells.Item(64, 4).Value = 44540
$ws.Cells.Item(64, 5).Value = 9
$ws.Cells.Item(64, 6).Value = 100112052
$ws.Cells.Item(64, 7).Value = "Albahaca"
$ws.Cells.Item(64, 8).Value = "Sin especificar"
$ws.Cells.Item(64, 9).Value = "Primera"
$ws.Cells.Item(64, 10).Value = 25
$ws.Cells.Item(64, 11).Value = 5000
$ws.Cells.Item(64, 12).Value = 5000
$ws.Cells.Item(64, 13).Value = 5000
$ws.Cells.Item(64, 14).Value = "$/paquete"
$ws.Cells.Item(64, 15).Value = "Región del Maule"
$ws.Cells.Item(64, 16).Value = 5000
$ws.Cells.Item(64, 17).Value = 1
$ws.Cells.Item(64, 18).Value = "Hortaliza"
